# Rename the column-A header of the loads worksheet from "Loads" to "Load Name".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Load Name"
